$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 20881.8
$ws.Range("I21").Value = 25545.428
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 25545.428
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = -25077.428
$ws.Range("N21").Value = -10936
$ws.Range("H23").Value = 20881.8
$ws.Range("I23").Value = 25545.428
$ws.Range("J23").Value = 10000
$ws.Range("K23").Value = 25545.428
$ws.Range("L23").Value = 10000
$ws.Range("M23").Value = -25311.428
$ws.Range("N23").Value = -10468
$ws.Range("H47").Value = 31500
$ws.Range("J47").Value = 48000
$ws.Range("L47").Value = 48000
$ws.Range("N47").Value = -49944
$ws.Range("H138").Value = 2415.3906
$ws.Range("I138").Value = 1248.25
$ws.Range("J138").Value = 2945.9092
$ws.Range("K138").Value = 3744.75
$ws.Range("L138").Value = 8837.7276
$ws.Range("M138").Value = 1395.25
$ws.Range("N138").Value = -19117.7276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13542.582
$ws.Range("I32").Value = 13978.441
$ws.Range("K32").Value = 13978.441
$ws.Range("M32").Value = -13691.441
$ws.Range("H61").Value = 2792.652
$ws.Range("I61").Value = 2262.2
$ws.Range("J61").Value = 3200.6924
$ws.Range("K61").Value = 2262.2
$ws.Range("L61").Value = 3200.6924
$ws.Range("M61").Value = -2050.2
$ws.Range("N61").Value = -3624.6924
$ws.Range("H88").Value = 1797003
$ws.Range("I88").Value = 8000
$ws.Range("J88").Value = 2020628.2
$ws.Range("K88").Value = 8000
$ws.Range("L88").Value = 2020628.2
$ws.Range("M88").Value = -7594
$ws.Range("N88").Value = -2021440.2
$ws.Range("H91").Value = 1797003
$ws.Range("I91").Value = 8000
$ws.Range("J91").Value = 2020628.2
$ws.Range("K91").Value = 8000
$ws.Range("L91").Value = 2020628.2
$ws.Range("M91").Value = -6596
$ws.Range("N91").Value = -2023436.2
$ws.Range("H102").Value = 27382.5
$ws.Range("J102").Value = 79600
$ws.Range("L102").Value = 79600
$ws.Range("N102").Value = -82844
$ws.Range("H136").Value = 2792.652
$ws.Range("I136").Value = 2262.2
$ws.Range("J136").Value = 3200.6924
$ws.Range("K136").Value = 6786.599999999999
$ws.Range("L136").Value = 9602.0772
$ws.Range("M136").Value = -4236.599999999999
$ws.Range("N136").Value = -14702.0772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 46780
$ws.Range("J81").Value = 46780
$ws.Range("L81").Value = 46780
$ws.Range("N81").Value = -48902
$ws.Range("H84").Value = 46780
$ws.Range("J84").Value = 46780
$ws.Range("L84").Value = 140340
$ws.Range("N84").Value = -150948
$ws.Range("H86").Value = 4841.6665
$ws.Range("I86").Value = 4066.1667
$ws.Range("J86").Value = 5617.1665
$ws.Range("K86").Value = 4066.1667
$ws.Range("L86").Value = 5617.1665
$ws.Range("M86").Value = -2943.1667
$ws.Range("N86").Value = -7863.1665
$ws.Range("H89").Value = 4841.6665
$ws.Range("I89").Value = 4066.1667
$ws.Range("J89").Value = 5617.1665
$ws.Range("K89").Value = 20330.8335
$ws.Range("L89").Value = 28085.8325
$ws.Range("M89").Value = -14714.8335
$ws.Range("N89").Value = -39317.8325
$ws.Range("H99").Value = 2114.8276
$ws.Range("I99").Value = 1978.0952
$ws.Range("J99").Value = 2473.75
$ws.Range("K99").Value = 1978.0952
$ws.Range("L99").Value = 2473.75
$ws.Range("M99").Value = -480.0952
$ws.Range("N99").Value = -5469.75
$ws.Range("H134").Value = 1920.4
$ws.Range("I134").Value = 1460.65
$ws.Range("J134").Value = 3759.4
$ws.Range("K134").Value = 4381.950000000001
$ws.Range("L134").Value = 11278.2
$ws.Range("M134").Value = -1846.950000000001
$ws.Range("N134").Value = -16348.2
$ws.Range("H138").Value = 48678.184
$ws.Range("J138").Value = 48678.184
$ws.Range("L138").Value = 48678.184
$ws.Range("N138").Value = -58958.184

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1405.069
$ws.Range("I68").Value = 1116.9
$ws.Range("J68").Value = 1556.7368
$ws.Range("K68").Value = 3350.7
$ws.Range("L68").Value = 4670.2104
$ws.Range("M68").Value = -2539.7
$ws.Range("N68").Value = -6292.2104
$ws.Range("H71").Value = 1405.069
$ws.Range("I71").Value = 1116.9
$ws.Range("J71").Value = 1556.7368
$ws.Range("K71").Value = 10052.1
$ws.Range("L71").Value = 14010.6312
$ws.Range("M71").Value = -5996.1
$ws.Range("N71").Value = -22122.6312
$ws.Range("H107").Value = 2747.1074
$ws.Range("I107").Value = 3852.0356
$ws.Range("J107").Value = 2271.1384
$ws.Range("K107").Value = 11556.1068
$ws.Range("L107").Value = 6813.415199999999
$ws.Range("M107").Value = -9636.106800000001
$ws.Range("N107").Value = -10653.4152
$ws.Range("H113").Value = 13725.625
$ws.Range("I113").Value = 14971.571
$ws.Range("K113").Value = 44914.713
$ws.Range("M113").Value = -42744.713
$ws.Range("H134").Value = 32361566
$ws.Range("I134").Value = 45597510
$ws.Range("J134").Value = 7029.1113
$ws.Range("K134").Value = 136792530
$ws.Range("L134").Value = 21087.3339
$ws.Range("M134").Value = -136787460
$ws.Range("N134").Value = -31227.3339
$ws.Range("H139").Value = 112778.3
$ws.Range("I139").Value = 112778.3
$ws.Range("K139").Value = 338334.9
$ws.Range("M139").Value = -333194.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5568.3184
$ws.Range("I70").Value = 5510.4
$ws.Range("K70").Value = 5510.4
$ws.Range("M70").Value = -5240.4
$ws.Range("H73").Value = 5568.3184
$ws.Range("I73").Value = 5510.4
$ws.Range("K73").Value = 5510.4
$ws.Range("M73").Value = -4574.4
$ws.Range("H80").Value = 9487.5
$ws.Range("I80").Value = 6975
$ws.Range("J80").Value = 12000
$ws.Range("K80").Value = 6975
$ws.Range("L80").Value = 12000
$ws.Range("M80").Value = -5977
$ws.Range("N80").Value = -13996
$ws.Range("H83").Value = 9487.5
$ws.Range("I83").Value = 6975
$ws.Range("J83").Value = 12000
$ws.Range("K83").Value = 34875
$ws.Range("L83").Value = 60000
$ws.Range("M83").Value = -29883
$ws.Range("N83").Value = -69984

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 585.3333
$ws.Range("I100").Value = 427.25
$ws.Range("J100").Value = 901.5
$ws.Range("K100").Value = 854.5
$ws.Range("L100").Value = 1803
$ws.Range("M100").Value = -313.5
$ws.Range("N100").Value = -2885
$ws.Range("H126").Value = 3678130.5
$ws.Range("I126").Value = 4903510.5
$ws.Range("K126").Value = 14710531.5
$ws.Range("M126").Value = -14708061.5
